$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2").Value = 613003
